# Reorders the "6-cylinder" summary rows of the mtcars table.
#
# Before:
#   row5 = cyl=4, vs=0                (style 11/12/12/13/14/14/14/9)
#   row6 = cyl=6, vs=1  (merged A6:A7/B6:B7/C6:C7, style 15/16/16/13/14/14/14/9)
#   row7 = cyl=6, vs=0  (A7 blank, only D:G filled)     <- the other half of the merge
#   row8 = cyl=8, vs=1
#
# After:
#   row5 = cyl=6, vs=0  (what used to be the "blank A" half of the merged block,
#                         now a standalone row re-using row6's numeric D:G values)
#   row6 = cyl=4, vs=1  (the old row5 data, plain style, B/C now =1)
#   row7 = cyl=6, vs=1  (the old row6 data, now unmerged / plain style)
#   row8 = cyl=8, vs=1  (unchanged)
#
# and the A6:A7 / B6:B7 / C6:C7 merges go away.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the "before" values we still need (rows get overwritten in place) ---
# NB: use .Value2 for reads -- .Value (getter, no-paren call) returns a
# property-descriptor string in this host rather than invoking the getter.
$d5 = $ws.Range("D5").Value2
$e5 = $ws.Range("E5").Value2
$f5 = $ws.Range("F5").Value2
$g5 = $ws.Range("G5").Value2

$d6 = $ws.Range("D6").Value2
$e6 = $ws.Range("E6").Value2
$f6 = $ws.Range("F6").Value2
$g6 = $ws.Range("G6").Value2

# --- unmerge the old 6-cylinder block (A6:A7, B6:B7, C6:C7) ---
$ws.Range("A6:A7").UnMerge()
$ws.Range("B6:B7").UnMerge()
$ws.Range("C6:C7").UnMerge()

# --- fix up styling so A6/B6/C6 use the plain (non-merged) formats, matching A5/B5/C5 ---
$ws.Range("A5").Copy($ws.Range("A6"))
$ws.Range("B5").Copy($ws.Range("B6"))
$ws.Range("C5").Copy($ws.Range("C6"))

# --- A7 was an empty cell (s="9"); give it the same plain style as A5/A8 ---
$ws.Range("A5").Copy($ws.Range("A7"))
# --- B7/C7 did not exist before; create them with the plain style too ---
$ws.Range("B5").Copy($ws.Range("B7"))
$ws.Range("C5").Copy($ws.Range("C7"))

# --- now write the new row 5: old row6's "cyl=6, vs=0" numbers ---
$ws.Range("A5").Value = 6
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = $d6
$ws.Range("E5").Value = $e6
$ws.Range("F5").Value = $f6
$ws.Range("G5").Value = $g6

# --- row 6: old row5's "cyl=4" numbers (am corrected to 1) ---
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = $d5
$ws.Range("E6").Value = $e5
$ws.Range("F6").Value = $f5
$ws.Range("G6").Value = $g5

# --- row 7: old row6's "cyl=6, vs=1" label columns; D:G (the 112.67 group) stay as-is ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 1
